$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) idPagamento (column D) was filled in for the existing rows 41-46 that
#    previously had a blank value.
# ---------------------------------------------------------------------------
# Keep these as TEXT (they are long numeric-looking ids), matching the other
# "idPagamento" cells already in the sheet.
$ws.Range("D41:D50").NumberFormat = "@"

$dPagamentos = @{
    41 = "77280262910"
    42 = "77282715892"
    43 = "77051216465"
    44 = "77298316438"
    45 = "77055139631"
    46 = "77298317798"
}
foreach ($r in $dPagamentos.Keys) {
    $ws.Cells.Item($r, 4).Value = $dPagamentos[$r]
}

# ---------------------------------------------------------------------------
# 2) Ten new rows (47-56) were appended for "Vitor Ito", each with the usual
#    N1..N10 sequence and a "Pagamento" = "Não" flag. Column C (Telefone)
#    holds numeric-looking text, so force it to stay TEXT as well.
# ---------------------------------------------------------------------------
$ws.Range("C47:C56").NumberFormat = "@"

$newRows = @(
    @{ Row = 47; C = "11988776655";  D = "77299240172"; Nums = 1..10 }
    @{ Row = 48; C = "11988776655";  D = "77056098875"; Nums = 1..10 }
    @{ Row = 49; C = "11977665544";  D = "77341583624"; Nums = 1..10 }
    @{ Row = 50; C = "11988776655";  D = "77342981440"; Nums = 1..10 }
    @{ Row = 51; C = "11988776655";  D = $null;          Nums = 1..10 }
    @{ Row = 52; C = "11988776655";  D = $null;          Nums = 1..10 }
    @{ Row = 53; C = "11988776655";  D = $null;          Nums = 1..10 }
    @{ Row = 54; C = "11988776655";  D = $null;          Nums = 11..20 }
    @{ Row = 55; C = "119988776655"; D = $null;          Nums = 1..10 }
    @{ Row = 56; C = "11977665544";  D = $null;          Nums = 1..10 }
)

foreach ($item in $newRows) {
    $r = $item.Row

    $ws.Cells.Item($r, 1).Value = "Vitor Ito"      # A - Nome
    $ws.Cells.Item($r, 2).Value = 1578424633        # B - ID
    $ws.Cells.Item($r, 3).Value = $item.C            # C - Telefone

    if ($null -ne $item.D) {
        $ws.Cells.Item($r, 4).Value = $item.D        # D - idPagamento
    }

    for ($i = 0; $i -lt 10; $i++) {
        $ws.Cells.Item($r, 5 + $i).Value = $item.Nums[$i]   # E..N - N1..N10
    }

    $ws.Cells.Item($r, 15).Value = "Não"             # O - Pagamento
}
